$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GOOG")

# Row 7 - Property, Plant, Equipment (Net)
$ws.Range("B7").Value = 100204000000.0
$ws.Range("G7").Value = 73646000000.0

# Row 8 - Long-Term Investments
$ws.Range("B8").Value = 25757000000.0

# Row 11 - Long-term assets (Other)
$ws.Range("B11").Value = 3704000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 3819000000.0
$ws.Range("G15").Value = 5561000000.0

# Row 17 - Current Part of Debt (was blank inline string, now a number)
$ws.Range("B17").Value = 2978000000.0

# Row 20 - Other current liabilities (was blank inline string, now a number)
$ws.Range("B20").Value = 45781000000.0

# Row 22 - Long Term Debt (Total)
$ws.Range("B22").Value = 25269000000.0
$ws.Range("G22").Value = 4554000000.0

# Row 24 - Long Term Tax Liability (Deferred)
$ws.Range("B24").Value = 4377000000.0
$ws.Range("G24").Value = 1701000000.0

# Row 29 - Common Stock (Net)
$ws.Range("B29").Value = 1000000.0
$ws.Range("G29").Value = 688000.0

# Row 36 - Net Debt
$ws.Range("B36").Value = -106857000000.0

# Row 37 - Total Debt
$ws.Range("B37").Value = 28247000000.0
